$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, bordered, centered/top aligned)
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data rows 2-61 get the team record values
for ($r = 2; $r -le 61; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
